$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values that changed.
# Values that look numeric are prefixed with a leading apostrophe so Excel
# stores/keeps them as text (matching the original text-formatted prices),
# rather than silently converting/rounding them to a number.
$ws.Range("D2").Value = "28.015.23"
$ws.Range("D3").Value = "1.826.53"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D5").Value = "'328.53"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D7").Value = "'0.4633"
$ws.Range("D8").Value = "'0.3861"
$ws.Range("D9").Value = "'46.08"
$ws.Range("D10").Value = "'0.07879"
$ws.Range("D11").Value = "'0.9598"
$ws.Range("D12").Value = "'21.86"
$ws.Range("D13").Value = "1.837.86"
$ws.Range("D14").Value = "'5.652"
$ws.Range("D15").Value = "'6.882"
$ws.Range("D16").Value = "'0.06842"
$ws.Range("D18").Value = "'86.42"
$ws.Range("D19").Value = "'0.000009964"
$ws.Range("D20").Value = "'16.64"
$ws.Range("D21").Value = "'1.001"
$ws.Range("D22").Value = "28.036.91"
$ws.Range("D23").Value = "'5.318"
$ws.Range("D24").Value = "'10.98"
$ws.Range("D25").Value = "'2.090"
$ws.Range("D26").Value = "2.052.98"
$ws.Range("D27").Value = "'152.17"
$ws.Range("D28").Value = "'19.14"
$ws.Range("D29").Value = "'5.744"
$ws.Range("D30").Value = "'1.968"
$ws.Range("D31").Value = "'116.64"
$ws.Range("D32").Value = "'0.9365"
$ws.Range("D33").Value = "'0.09219"
$ws.Range("D35").Value = "'1.317"
$ws.Range("D36").Value = "'3.344"
$ws.Range("D37").Value = "'0.05929"
$ws.Range("D39").Value = "'1.144"
$ws.Range("D40").Value = "'0.9996"
$ws.Range("D41").Value = "'7.609"
$ws.Range("D42").Value = "'0.5572"
$ws.Range("D43").Value = "'9.904"
$ws.Range("D44").Value = "'0.1763"
$ws.Range("D45").Value = "'1.225"
$ws.Range("D46").Value = "'2.221"
$ws.Range("D47").Value = "'11.60"
$ws.Range("D48").Value = "'0.5258"
$ws.Range("D49").Value = "'0.07007"
$ws.Range("D51").Value = "'111.26"

# Update Volume(1h) (column E) values, preserving the two leading/trailing spaces
$ws.Range("E2").Value = "  -5.12%  "
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("E5").Value = "  -3.04%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("E12").Value = "  -5.72%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  -5.04%  "
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("E24").Value = "  -5.48%  "
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("E29").Value = "  -11.17%  "
$ws.Range("E30").Value = "  -4.44%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("E36").Value = "  -5.27%  "
$ws.Range("E37").Value = "  -7.22%  "
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("E43").Value = "  -6.35%  "
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("E46").Value = "  -8.33%  "
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("E50").Value = "  -6.86%  "
$ws.Range("E51").Value = "  -4.42%  "
